# Fruta / hortaliza, semanal
# Insert a new weekly price block (3 rows: Especial/Primera/Segunda) for
# Femacal de La Calera - Kiwi, dated 2023-07-25 (serial 45132), region
# "Provincia de Curicó", right before the existing row 1092 data block.
# All subsequent rows shift down by 3 (sheet dimension grows from
# A1:T1201 to A1:T1204).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows above row 1092, pushing the rest of the table down.
$ws.Range("A1092:T1094").EntireRow.Insert()

# Populate the 3 newly inserted rows with the new weekly price data.

# Row 1092: Calidad "Especial"
$ws.Range("A1092").Value = 3
$ws.Range("B1092").Value = "Femacal de La Calera"
$ws.Range("C1092").Value = "Coquimbo"
$ws.Range("D1092").Value = 45132
$ws.Range("E1092").Value = 5
$ws.Range("F1092").Value = "Fruta"
$ws.Range("G1092").Value = 100101
$ws.Range("H1092").Value = "Berries"
$ws.Range("I1092").Value = 100101007
$ws.Range("J1092").Value = "Kiwi"
$ws.Range("K1092").Value = "Hayward"
$ws.Range("L1092").Value = "Especial"
$ws.Range("M1092").Value = 75
$ws.Range("N1092").Value = 8000
$ws.Range("O1092").Value = 8000
$ws.Range("P1092").Value = 8000
$ws.Range("Q1092").Value = "$/bandeja 10 kilos"
$ws.Range("R1092").Value = "Provincia de Curicó"
$ws.Range("S1092").Value = 800
$ws.Range("T1092").Value = 10

# Row 1093: Calidad "Primera"
$ws.Range("A1093").Value = 3
$ws.Range("B1093").Value = "Femacal de La Calera"
$ws.Range("C1093").Value = "Coquimbo"
$ws.Range("D1093").Value = 45132
$ws.Range("E1093").Value = 5
$ws.Range("F1093").Value = "Fruta"
$ws.Range("G1093").Value = 100101
$ws.Range("H1093").Value = "Berries"
$ws.Range("I1093").Value = 100101007
$ws.Range("J1093").Value = "Kiwi"
$ws.Range("K1093").Value = "Hayward"
$ws.Range("L1093").Value = "Primera"
$ws.Range("M1093").Value = 68
$ws.Range("N1093").Value = 7000
$ws.Range("O1093").Value = 7000
$ws.Range("P1093").Value = 7000
$ws.Range("Q1093").Value = "$/bandeja 10 kilos"
$ws.Range("R1093").Value = "Provincia de Curicó"
$ws.Range("S1093").Value = 700
$ws.Range("T1093").Value = 10

# Row 1094: Calidad "Segunda"
$ws.Range("A1094").Value = 3
$ws.Range("B1094").Value = "Femacal de La Calera"
$ws.Range("C1094").Value = "Coquimbo"
$ws.Range("D1094").Value = 45132
$ws.Range("E1094").Value = 5
$ws.Range("F1094").Value = "Fruta"
$ws.Range("G1094").Value = 100101
$ws.Range("H1094").Value = "Berries"
$ws.Range("I1094").Value = 100101007
$ws.Range("J1094").Value = "Kiwi"
$ws.Range("K1094").Value = "Hayward"
$ws.Range("L1094").Value = "Segunda"
$ws.Range("M1094").Value = 65
$ws.Range("N1094").Value = 6000
$ws.Range("O1094").Value = 6000
$ws.Range("P1094").Value = 6000
$ws.Range("Q1094").Value = "$/bandeja 10 kilos"
$ws.Range("R1094").Value = "Provincia de Curicó"
$ws.Range("S1094").Value = 600
$ws.Range("T1094").Value = 10
